$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProductLoanInput")
$ws.Range("B17").Value = "Penalties, Fees, Interest, Principal order"
$ws.Range("B17").HorizontalAlignment = -4131
$ws.Range("B17").VerticalAlignment = -4160
$ws.Activate()
[void]$ws.Range("B17").Select()
